$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Severity / Priority" header to "Severity" (also updates the
# shared string table and the linked table column name automatically).
$ws.Range("H2").Value = "Severity"

# The header row grew taller after the edit (wrapped text needs more room).
$ws.Rows(2).RowHeight = 45

# Update the active selection left behind by the editor.
$ws.Range("I3").Select()
